$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.311.52'
$ws.Range("E2").Value = '  -3.18%  '
$ws.Range("D3").Value = '3.299.48'
$ws.Range("E3").Value = '  -3.81%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '556.72'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.88%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '141.62'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -7.16%  '
$ws.Range("E7").Value = '  -0.15%  '
$ws.Range("D8").Value = '3.303.63'
$ws.Range("E8").Value = '  -3.71%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.467'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.27%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.84'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.60%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.118'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -5.03%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.408'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.36%  '
$ws.Range("D13").Value = '3.871.81'
$ws.Range("E13").Value = '  -3.51%  '
$ws.Range("E14").Value = '  +0.54%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '26.98'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -5.85%  '
$ws.Range("D16").Value = '3.300.04'
$ws.Range("E16").Value = '  -4.84%  '
$ws.Range("E17").Value = '  -4.14%  '
$ws.Range("D18").Value = '60.353.23'
$ws.Range("E18").Value = '  -3.17%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.11'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -6.35%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.97'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.98%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '8.57'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.61%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '375.12'
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '74.36'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.17%  '
$ws.Range("E24").Value = '  +0.09%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.533'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -6.71%  '
$ws.Range("D26").Value = '3.438.97'
$ws.Range("E26").Value = '  -3.42%  '
$ws.Range("E27").Value = '  -9.79%  '
$ws.Range("E28").Value = '  -4.65%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.27%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.13'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -7.42%  '
$ws.Range("E31").Value = '  -0.10%  '
$ws.Range("E32").Value = '  -4.41%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.53'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -5.64%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '22.60'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.79%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.23'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -8.42%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.10'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -6.94%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '166.92'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.20%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.52'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -6.18%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.65'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -4.16%  '
$ws.Range("D40").Value = '3.336.77'
$ws.Range("E40").Value = '  -3.64%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '26.55'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -14.62%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0728'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -7.52%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '41.91'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.14%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.750'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.90%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.11'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -6.73%  '
$ws.Range("E46").Value = '  -6.32%  '
$ws.Range("E47").Value = '  -7.18%  '
$ws.Range("B48").Value = 'FirstDigitalUSD'
$ws.Range("C48").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.00'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.00%  '
$ws.Range("B49").Value = 'Maker'
$ws.Range("C49").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D49").Value = '2.345.49'
$ws.Range("E49").Value = '  -7.72%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.37'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -7.59%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '21.24'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -6.14%  '
